# Auto-generated Excel COM-interop script applying the "Goblin Profits" value updates
# described by the commit diff. Each row touches only columns H-N (all numeric,
# no formulas in this workbook), so we set/clear cells directly by A1 reference.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 874.8182
$ws.Range("J19").Value = 725.4286
$ws.Range("L19").Value = 725.4286
$ws.Range("N19").Value = -1075.4286

$ws.Range("H20").Value = 1899.6666
$ws.Range("I20").Value = 1899.6666
$ws.Range("K20").Value = 1899.6666
$ws.Range("M20").Value = -1669.6666

$ws.Range("H28").Value = 6489.4375
$ws.Range("J28").Value = 332
$ws.Range("L28").Value = 332
$ws.Range("N28").Value = -1302

$ws.Range("H35").Value = 1899.6666
$ws.Range("I35").Value = 1899.6666
$ws.Range("K35").Value = 1899.6666
$ws.Range("M35").Value = -1520.6666

$ws.Range("H40").Value = 3331
$ws.Range("I40").Value = 2595.8
$ws.Range("J40").Value = 4250
$ws.Range("K40").Value = 2595.8
$ws.Range("L40").Value = 4250
$ws.Range("M40").Value = -2420.8
$ws.Range("N40").Value = -4600

$ws.Range("H86").Value = 2082.9092
$ws.Range("I86").Value = 2099.4666
$ws.Range("K86").Value = 2099.4666
$ws.Range("M86").Value = -976.4666000000002

$ws.Range("H89").Value = 2082.9092
$ws.Range("I89").Value = 2099.4666
$ws.Range("K89").Value = 10497.333
$ws.Range("M89").Value = -4881.333000000001

$ws.Range("H98").Value = 6700.5557
$ws.Range("I98").Value = 9116.684999999999
$ws.Range("K98").Value = 9116.684999999999
$ws.Range("M98").Value = -7618.684999999999

$ws.Range("H122").Value = 6700.5557
$ws.Range("I122").Value = 9116.684999999999
$ws.Range("K122").Value = 27350.055
$ws.Range("M122").Value = -24900.055

$ws.Range("H135").Value = 2000
$ws.Range("I135").Value = 2000
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18000
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("M135").Value = -15465

$ws.Range("H138").Value = 1750.7462
$ws.Range("I138").Value = 1232.1428
$ws.Range("J138").Value = 1887.7358
$ws.Range("K138").Value = 3696.4284
$ws.Range("L138").Value = 5663.207399999999
$ws.Range("M138").Value = 1443.5716
$ws.Range("N138").Value = -15943.2074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3332.6667
$ws.Range("I3").Value = 3332.6667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3332.6667
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("M3").Value = -3217.6667

$ws.Range("H5").Value = 192.91667
$ws.Range("I5").Value = 66.42856999999999
$ws.Range("J5").Value = 370
$ws.Range("K5").Value = 66.42856999999999
$ws.Range("L5").Value = 370
$ws.Range("M5").Value = 45.57143000000001
$ws.Range("N5").Value = -594

$ws.Range("H74").Value = 2494.0557
$ws.Range("J74").Value = 1496.2
$ws.Range("L74").Value = 1496.2
$ws.Range("N74").Value = -3244.2

$ws.Range("H77").Value = 2494.0557
$ws.Range("J77").Value = 1496.2
$ws.Range("L77").Value = 7481
$ws.Range("N77").Value = -16217

$ws.Range("H122").Value = 2443.0715
$ws.Range("I122").Value = 2299.625
$ws.Range("K122").Value = 6898.875
$ws.Range("M122").Value = -4448.875

$ws.Range("H132").Value = 7514.3887
$ws.Range("I132").Value = 8770.666999999999
$ws.Range("J132").Value = 1233
$ws.Range("K132").Value = 26312.001
$ws.Range("L132").Value = 3699
$ws.Range("M132").Value = -23782.001
$ws.Range("N132").Value = -8759

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 192.91667
$ws.Range("I4").Value = 66.42856999999999
$ws.Range("J4").Value = 370
$ws.Range("K4").Value = 66.42856999999999
$ws.Range("L4").Value = 370
$ws.Range("M4").Value = 48.57143000000001
$ws.Range("N4").Value = -600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 751.7778
$ws.Range("I16").Value = 751.7778
$ws.Range("K16").Value = 751.7778
$ws.Range("M16").Value = -464.7778

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("L26").Value = 0

$ws.Range("H31").Value = 2883.2307
$ws.Range("I31").Value = 1018.04
$ws.Range("J31").Value = 6213.9287
$ws.Range("K31").Value = 1018.04
$ws.Range("L31").Value = 6213.9287
$ws.Range("M31").Value = -723.04
$ws.Range("N31").Value = -6803.9287

$ws.Range("H34").Value = 2883.2307
$ws.Range("I34").Value = 1018.04
$ws.Range("J34").Value = 6213.9287
$ws.Range("K34").Value = 1018.04
$ws.Range("L34").Value = 6213.9287
$ws.Range("M34").Value = -816.04
$ws.Range("N34").Value = -6617.9287

$ws.Range("H51").Value = 12000
$ws.Range("I51").Value = 12000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 12000
$ws.Range("N51").ClearContents()
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -11264

$ws.Range("H58").Value = 1297.0588
$ws.Range("I58").Value = 1186.7273
$ws.Range("K58").Value = 1186.7273
$ws.Range("M58").Value = -983.7273

$ws.Range("H61").Value = 12000
$ws.Range("I61").Value = 12000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 12000
$ws.Range("N61").ClearContents()
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -11652

$ws.Range("H94").Value = 2230.8667
$ws.Range("J94").Value = 2254.25
$ws.Range("L94").Value = 2254.25
$ws.Range("N94").Value = -3156.25

$ws.Range("H113").Value = 751.7778
$ws.Range("I113").Value = 751.7778
$ws.Range("K113").Value = 751.7778
$ws.Range("M113").Value = 1418.2222

$ws.Range("H132").Value = 2339.2727
$ws.Range("I132").Value = 2289.1667
$ws.Range("K132").Value = 6867.500100000001
$ws.Range("M132").Value = -4337.500100000001

$ws.Range("H136").Value = 1297.0588
$ws.Range("I136").Value = 1186.7273
$ws.Range("K136").Value = 3560.1819
$ws.Range("M136").Value = -1010.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 145.64285
$ws.Range("J33").Value = 198.66667
$ws.Range("L33").Value = 1192.00002
$ws.Range("N33").Value = -1758.00002

$ws.Range("H46").Value = 11213198
$ws.Range("I46").Value = 13939516
$ws.Range("K46").Value = 41818548
$ws.Range("M46").Value = -41818457

$ws.Range("H49").Value = 642.2
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H99").Value = 15845
$ws.Range("I99").Value = 8284
$ws.Range("J99").Value = 24666.166
$ws.Range("K99").Value = 24852
$ws.Range("L99").Value = 73998.49800000001
$ws.Range("M99").Value = -22606
$ws.Range("N99").Value = -78490.49800000001

$ws.Range("H121").Value = 686.375
$ws.Range("I121").Value = 346.375
$ws.Range("K121").Value = 1039.125
$ws.Range("M121").Value = 270.875

$ws.Range("H129").Value = 3424.818
$ws.Range("J129").Value = 5966.6665
$ws.Range("L129").Value = 17899.9995
$ws.Range("N129").Value = -27899.9995

$ws.Range("H132").Value = 2964.9
$ws.Range("I132").Value = 2824.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 25420.5
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -22890.5
$ws.Range("N132").Value = -32060

$ws.Range("H137").Value = 4453.533
$ws.Range("J137").Value = 7021.125
$ws.Range("L137").Value = 21063.375
$ws.Range("N137").Value = -31263.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6126.154
$ws.Range("I80").Value = 4233
$ws.Range("J80").Value = 6694.1
$ws.Range("K80").Value = 4233
$ws.Range("L80").Value = 6694.1
$ws.Range("M80").Value = -3235
$ws.Range("N80").Value = -8690.1

$ws.Range("H83").Value = 6126.154
$ws.Range("I83").Value = 4233
$ws.Range("J83").Value = 6694.1
$ws.Range("K83").Value = 21165
$ws.Range("L83").Value = 33470.5
$ws.Range("M83").Value = -16173
$ws.Range("N83").Value = -43454.5

$ws.Range("H104").Value = 54987.5
$ws.Range("J104").Value = 54987.5
$ws.Range("L104").Value = 54987.5
$ws.Range("N104").Value = -61975.5

$ws.Range("H122").Value = 28388.223
$ws.Range("I122").Value = 31999.143
$ws.Range("J122").Value = 15750
$ws.Range("K122").Value = 95997.429
$ws.Range("L122").Value = 47250
$ws.Range("M122").Value = -93547.429
$ws.Range("N122").Value = -52150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3794.8572
$ws.Range("J46").Value = 4046.1
$ws.Range("L46").Value = 4046.1
$ws.Range("N46").Value = -4422.1

$ws.Range("H93").Value = 6826.8184
$ws.Range("I93").Value = 5619
$ws.Range("J93").Value = 7833.3335
$ws.Range("K93").Value = 5619
$ws.Range("L93").Value = 7833.3335
$ws.Range("M93").Value = -4371
$ws.Range("N93").Value = -10329.3335

$ws.Range("H103").Value = 24999.5
$ws.Range("J103").Value = 24999.5
$ws.Range("L103").Value = 24999.5
$ws.Range("N103").Value = -27343.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 22499.75
$ws.Range("J33").Value = 22499.75
$ws.Range("L33").Value = 22499.75
$ws.Range("N33").Value = -22999.75

$ws.Range("H36").Value = 22499.75
$ws.Range("J36").Value = 22499.75
$ws.Range("L36").Value = 22499.75
$ws.Range("N36").Value = -22999.75

$ws.Range("H126").Value = 1945.091
$ws.Range("I126").Value = 1624.5
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 4873.5
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -2403.5
$ws.Range("N126").Value = -13340

$ws.Range("H132").Value = 10537.52
$ws.Range("I132").Value = 16229.934
$ws.Range("J132").Value = 1998.9
$ws.Range("K132").Value = 48689.802
$ws.Range("L132").Value = 5996.700000000001
$ws.Range("M132").Value = -46159.802
$ws.Range("N132").Value = -11056.7
